$d = $word.ActiveDocument

# --- 1) Merge "(True/False) " + "The regression equation for adding potential
#         confounders ..." into a single run ---
$d.Content.Find.Execute(
    "(True/False) The regression equation for adding potential confounders is the same as the regression equation for adding potential effect modifiers.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(True/False) The regression equation for adding potential confounders is the same as the regression equation for adding potential effect modifiers.",
    2) | Out-Null

# --- 2) Merge "(True/False) " + "Given an outcome Y and covariates X and Z, ..." ---
$d.Content.Find.Execute(
    "(True/False) Given an outcome Y and covariates X and Z, if the linearity assumptions ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(True/False) Given an outcome Y and covariates X and Z, if the linearity assumptions ",
    2) | Out-Null

# --- 3) Merge "(Multiple choice) " + "Taylor is studying the association ..." ---
$d.Content.Find.Execute(
    "(Multiple choice) Taylor is studying the association between continuous quantitative variables X (predictor of interest) and Y (outcome), ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Multiple choice) Taylor is studying the association between continuous quantitative variables X (predictor of interest) and Y (outcome), ",
    2) | Out-Null

# --- 4) Merge "A" + "mong the following, which is the most reasonable ..." ---
$d.Content.Find.Execute(
    "Among the following, which is the most reasonable conclusion to draw from this plot?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Among the following, which is the most reasonable conclusion to draw from this plot?",
    2) | Out-Null

# --- 5) & 6) Mark the two embedded-picture runs as NoProof ($w:noProof/$) ---
$shapes = $d.InlineShapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shapes.Item($i).Range.NoProofing = -1
}

# --- 7) Insert a new paragraph (red note) right before the "(Short answer)
#         Taylor is building a prediction model ..." question, and merge its
#         own split runs in the same step ---
$paras = $d.Paragraphs
$targetIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Taylor is building a prediction model*") {
        $targetIndex = $i
        break
    }
}

$anchor = $paras.Item($targetIndex - 1)
$anchor.Range.InsertParagraphAfter() | Out-Null

$paras2 = $d.Paragraphs
$newPara = $paras2.Item($targetIndex)
$newPara.Range.Text = "(Everyone got credit for this question because it was confusingly worded)"
$newPara.Format.LeftIndent = 36
$newPara.Range.Font.Color = 255

# --- 8) Merge "(Short answer) " + "Taylor is building a prediction model ..." ---
$d.Content.Find.Execute(
    "(Short answer) Taylor is building a prediction model to predict the number of hours her cat Alice will sleep on a given night. In her model, she includes the following predictors:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(Short answer) Taylor is building a prediction model to predict the number of hours her cat Alice will sleep on a given night. In her model, she includes the following predictors:",
    2) | Out-Null
